$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.141.28"
$ws.Range("E2").Value = "  +5.40%  "
$ws.Range("D3").Value = "1.779.32"
$ws.Range("E3").Value = "  +2.89%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.69"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4921"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2671"
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06264"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "1.779.18"
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07027"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6277"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.649"
$ws.Range("E14").Value = "  +3.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "79.84"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").Value = "28.125.15"
$ws.Range("E16").Value = "  +6.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007242"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +5.51%  "
$ws.Range("D21").Value = "2.007.10"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.565"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.737"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.229"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.42"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.860"
$ws.Range("E27").Value = "  +4.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.31"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.386"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.175"
$ws.Range("E30").Value = "  +6.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08257"
$ws.Range("E31").Value = "  +3.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.758"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04879"
$ws.Range("E33").Value = "  +9.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.071"
$ws.Range("E34").Value = "  +7.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.617"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6511"
$ws.Range("E36").Value = "  +4.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9490"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.596"
$ws.Range("E38").Value = "  +7.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.040"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.892"
$ws.Range("E40").Value = "  +4.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01549"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.76"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3981"
$ws.Range("E44").Value = "  +3.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.175"
$ws.Range("E45").Value = "  +4.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1210"
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05434"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.988"
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.297"
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.71"
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.86"
$ws.Range("E51").Value = "  +2.20%  "
